$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly price record needs to be inserted at row 23, which pushes the
# existing rows 23-36 down to 24-37 (the sheet grows from A1:R36 to A1:R37).
$ws.Rows(23).Insert()

# Populate the newly inserted row 23 with the new "Locoto" price observation.
$ws.Range("A23").Value = 10
$ws.Range("B23").Value = "Vega Modelo de Temuco"
$ws.Range("C23").Value = "La Araucanía"
$ws.Range("D23").Value2 = 44784
$ws.Range("E23").Value = 9
$ws.Range("F23").Value = 100112042
$ws.Range("G23").Value = "Locoto"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 180
$ws.Range("K23").Value = 2700
$ws.Range("L23").Value = 2700
$ws.Range("M23").Value = 2700
$ws.Range("N23").Value = '$/kilo'
$ws.Range("O23").Value = "Región de Arica y Parinacota"
$ws.Range("P23").Value = 2700
$ws.Range("Q23").Value = 1
$ws.Range("R23").Value = "Hortaliza"
